$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# Row 24: 40. Combination Sum 2 (Backtracking)
# ---------------------------------------------------------------------------
$row24 = $lo.ListRows.Add()
$r24 = $row24.Range

$r24.Cells.Item(1, 1).Value = "40. Combination Sum 2"

$r24.Cells.Item(1, 2).Value = "Medium"
$ws.Range("B2").Copy()
$r24.Cells.Item(1, 2).PasteSpecial(-4122)

$r24.Cells.Item(1, 3).Value = "Backtracking"

$r24.Cells.Item(1, 4).Value = "Use distance from target as a loop condition. Sort the array first, and consider prev to handle duplicates."

$r24.Cells.Item(1, 5).Value = "https://leetcode.com/problems/combination-sum-ii/solutions/16878/combination-sum-i-ii-and-iii-java-solution-see-the-similarities-yourself/ "
$ws.Hyperlinks.Add($r24.Cells.Item(1, 5), "https://leetcode.com/problems/combination-sum-ii/solutions/16878/combination-sum-i-ii-and-iii-java-solution-see-the-similarities-yourself/") | Out-Null
$r24.Cells.Item(1, 5).Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Row 25: 347. Top K Frequent Elements (Arrays)
# ---------------------------------------------------------------------------
$row25 = $lo.ListRows.Add()
$r25 = $row25.Range

$r25.Cells.Item(1, 1).Value = "347. Top K Frequent Elements"

$r25.Cells.Item(1, 2).Value = "Medium"
$ws.Range("B2").Copy()
$r25.Cells.Item(1, 2).PasteSpecial(-4122)

$r25.Cells.Item(1, 3).Value = "Arrays"

$r25.Cells.Item(1, 4).Value = "Bucket sort. Use  counts as indices and numbers as values, and a hashmap to count occurrences."

$r25.Cells.Item(1, 5).Value = "https://leetcode.com/problems/top-k-frequent-elements/solutions/81602/java-o-n-solution-bucket-sort/ "
$ws.Hyperlinks.Add($r25.Cells.Item(1, 5), "https://leetcode.com/problems/top-k-frequent-elements/solutions/81602/java-o-n-solution-bucket-sort/ ") | Out-Null
$r25.Cells.Item(1, 5).Style = "Hyperlink"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Move the active selection down the way Excel would after typing through
# the new rows.
# ---------------------------------------------------------------------------
$ws.Range("D30").Select() | Out-Null
